# login_details.xlsx rework:
#  - rename Sheet1 -> Login_Details
#  - select the whole sheet (matches the recorded "select all" UI state)
#  - best-fit the three data columns (A/B/C) to their content widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Login_Details"

# "Select All" (Ctrl+A) so the saved view reflects the whole-sheet selection.
[void]$ws.Cells.Select()

# Column widths equivalent to Excel's AutoFit/"best fit" for this data:
#   A ~ 11.33 chars, B ~ 13.33 chars, C ~ 23.16 chars
$ws.Columns.Item(1).ColumnWidth = 10.498697916666666
$ws.Columns.Item(2).ColumnWidth = 12.498697916666666
$ws.Columns.Item(3).ColumnWidth = 22.330729166666668
